{"js": "// Remove the \"Input controls move to next row...\" bullet paragraph entirely\n// (fix formatting problem with field-with-errors class).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"Input controls move to next row when there is an error on the form (with FF and explorer only, chrome is fine)\";\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.trim() === target) {\n    paragraph.delete();\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Input controls move to next row...\" bullet paragraph entirely\n# (fix formatting problem with field-with-errors class).\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"Input controls move to next row when there is an error on the form (with FF and explorer only, chrome is fine)\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    $para = $range.Paragraphs(1)\n    $para.Range.Delete()\n}\n"}
